$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update frozen-pane view / selection state ---
$ws.Range("B458").Select()
$ws.Application.ActiveWindow.ScrollRow = 458
$ws.Range("S454").Select()

# --- Row 481 ---
$ws.Range("C481").Value = 9

# --- Row 482 ---
$ws.Range("C482").Value = 9

# --- Row 483 ---
$ws.Range("C483").Value = 8
$ws.Range("E483").Value = 3
$ws.Range("F483").Value = 1
$ws.Range("G483").Value = 2
$ws.Range("H483").Formula = "=IF(TODAY()>A482,G483+E483,"""")"
$ws.Range("L483").Value = 0
$ws.Range("M483").Value = 0

# --- Row 484 ---
$ws.Range("C484").Value = 6
$ws.Range("E484").Value = 2
$ws.Range("F484").Value = 1
$ws.Range("G484").Value = 4
$ws.Range("H484").Formula = "=IF(TODAY()>A483,G484+E484,"""")"
$ws.Range("L484").Value = 0
$ws.Range("M484").Value = 0

# --- Row 485 ---
$ws.Range("C485").Value = 0
$ws.Range("E485").Value = 1
$ws.Range("F485").Value = 1
$ws.Range("G485").Value = 5
$ws.Range("H485").Formula = "=IF(TODAY()>A484,G485+E485,"""")"
$ws.Range("L485").Value = 0
$ws.Range("M485").Value = 0
